# Disponibilidad.xlsx - "Actualizar 02-06-2021 00-34-35" automatic update.
# Appends one fresh availability-check pass (14 rows, one per monitored
# service) below the existing data, mirroring the same Nombre/URL/
# Disponibilidad pattern used by every previous pass, stamped with the new
# run's timestamp and carrying its own per-cell hyperlinks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = 44233.02387112618
$lastExistingRow = 897
$firstNewRow = 898

$newRows = @(
    @{ Row = 898;  Name = "Odoo";              Url = "https://www.dataintelligence-group.com/";                       SubAddress = $null },
    @{ Row = 899;  Name = "Blackbox";          Url = "https://serviciodashboard.azurewebsites.net/";                  SubAddress = $null },
    @{ Row = 900;  Name = "PowerBI";           Url = "https://powerbi.microsoft.com/es-es/";                          SubAddress = $null },
    @{ Row = 901;  Name = "Dropbox";           Url = "https://www.dropbox.com/";                                      SubAddress = $null },
    @{ Row = 902;  Name = "Odoo";              Url = "https://dataintelligence.store/";                               SubAddress = $null },
    @{ Row = 903;  Name = "GEE";               Url = "https://app-data-i.users.earthengine.app/";                     SubAddress = $null },
    @{ Row = 904;  Name = "UtilidadesOdoo";    Url = "https://odooutil.azurewebsites.net/";                           SubAddress = $null },
    @{ Row = 905;  Name = "Filtros Dashboard"; Url = "https://filtradordashboard.azurewebsites.net/";                 SubAddress = $null },
    @{ Row = 906;  Name = "MapStore";          Url = "https://ide.dataintelligence-group.com/mapstore/";              SubAddress = "/" },
    @{ Row = 907;  Name = "GeoServer";         Url = "https://ide.dataintelligence-group.com/geoserver/web/?0";       SubAddress = $null },
    @{ Row = 908;  Name = "Tomcat";            Url = "https://ide.dataintelligence-group.com/";                       SubAddress = $null },
    @{ Row = 909;  Name = "Shiny";             Url = "https://rpubs.com/dataintelligence/";                           SubAddress = $null },
    @{ Row = 910;  Name = "Github";            Url = "https://github.com/Sud-Austral/";                               SubAddress = $null },
    @{ Row = 911;  Name = "EZ Exporter";       Url = "https://ezexporter.highviewapps.com/exports/export-profile/";  SubAddress = $null }
)

foreach ($item in $newRows) {
    $r = $item.Row

    # Clone formatting (styles, number format, hyperlink font) from the row
    # directly above by copying the whole A:D block down one row at a time.
    $srcRow = $r - 1
    $ws.Range("A$srcRow`:D$srcRow").Copy($ws.Range("A$r`:D$r"))

    $ws.Range("A$r").Value = $item.Name
    $ws.Range("C$r").Value = "Disponible"
    $ws.Range("D$r").Value = $newTimestamp

    # Register the hyperlink relationship for column B, then restore the
    # display text and the shared "Hyperlink" cell style (Hyperlinks.Add
    # stamps its own variant xf, so re-applying the named style keeps it on
    # the same style used throughout the rest of the column).
    if ($item.SubAddress) {
        $ws.Hyperlinks.Add($ws.Range("B$r"), $item.Url, $item.SubAddress)
        $ws.Range("B$r").Value = $item.Url + "#" + $item.SubAddress
    } else {
        $ws.Hyperlinks.Add($ws.Range("B$r"), $item.Url)
        $ws.Range("B$r").Value = $item.Url
    }
    $ws.Range("B$r").Style = "Hyperlink"
}

Write-Output "Added rows $firstNewRow to $($newRows[-1].Row)"
